# diagnostic.xlsx: add a small "disconnected_elements" flag table
#   B1 = 0            (boxed / bold / centered header style)
#   A2 = 0            (same boxed style)
#   B2 = "disconnected_elements" (plain, shared string)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

# Build the boxed/bold/centered style once on B1 ...
$ws.Range("B1").Font.Bold = $true
$ws.Range("B1").HorizontalAlignment = -4108   # xlCenter
$ws.Range("B1").VerticalAlignment = -4160     # xlTop
$ws.Range("B1").Borders.LineStyle = 1         # xlContinuous
$ws.Range("B1").Borders.Weight = 2            # xlThin
$ws.Range("B1").Borders.ColorIndex = -4105    # xlColorIndexAutomatic

# ... then copy that exact style onto A2 so both cells share one cellXf
# instead of the engine minting a second, near-identical one.
$ws.Range("B1").Copy()
$ws.Range("A2").PasteSpecial(-4122)           # xlPasteFormats
$excel.CutCopyMode = 0
